$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 535
$ws.Range("I2").Value = 1303
$ws.Range("J2").Value = 5404
$ws.Range("K2").Value = 23
$ws.Range("L2").Value = 1436
$ws.Range("M2").Value = 102
$ws.Range("N2").Value = 990
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 27
$ws.Range("Q2").Value = 5
$ws.Range("R2").Value = 70
$ws.Range("S2").Value = 580
$ws.Range("T2").Value = 975
$ws.Range("U2").Value = 89
$ws.Range("V2").Value = 8566
$ws.Range("X2").Value = 8554
$ws.Range("Y2").Value = 10
$ws.Range("Z2").Value = 138
$ws.Range("AA2").Value = 44
